$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.568826
$ws.Range("H2").Value = 13.706478
$ws.Range("I2").Value = 0.08414374058887295
$ws.Range("J2").Value = 0.08414374058887293
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 9.261537239932
$ws.Range("R2").Value = 83.353835159388
$ws.Range("S2").Value = 0.0005550360573229706
$ws.Range("T2").Value = 0.0005550360573229703
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.568826
$ws.Range("H3").Value = 13.706478
$ws.Range("I3").Value = 0.08414374058887295
$ws.Range("J3").Value = 0.08414374058887293
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 1171.649409758872
$ws.Range("R3").Value = 10544.84468782985
$ws.Range("S3").Value = 0.07021595358419402
$ws.Range("T3").Value = 0.07021595358419401
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.568826
$ws.Range("H4").Value = 13.706478
$ws.Range("I4").Value = 0.08414374058887295
$ws.Range("J4").Value = 0.08414374058887293
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 223.1426756247741
$ws.Range("R4").Value = 2008.284080622966
$ws.Range("S4").Value = 0.01337275094735596
$ws.Range("T4").Value = 0.01337275094735595
$ws.Range("I5").Value = 0.6625544448906389
$ws.Range("J5").Value = 0.6625544448906389
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 72.92607414280533
$ws.Range("R5").Value = 656.334667285248
$ws.Range("S5").Value = 0.00437039765858162
$ws.Range("T5").Value = 0.00437039765858162
$ws.Range("I6").Value = 0.6625544448906389
$ws.Range("J6").Value = 0.6625544448906389
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("S6").Value = 0.5528859523461039
$ws.Range("T6").Value = 0.5528859523461038
$ws.Range("I7").Value = 0.6625544448906389
$ws.Range("J7").Value = 0.6625544448906389
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 1757.043014076971
$ws.Range("R7").Value = 15813.38712669274
$ws.Range("S7").Value = 0.1052980948859534
$ws.Range("T7").Value = 0.1052980948859533
$ws.Range("G8").Value = 13.75374933333333
$ws.Range("H8").Value = 41.261248
$ws.Range("I8").Value = 0.2533018145204882
$ws.Range("J8").Value = 0.2533018145204882
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 27.88043616442311
$ws.Range("R8").Value = 250.923925479808
$ws.Range("S8").Value = 0.001670850849514026
$ws.Range("T8").Value = 0.001670850849514025
$ws.Range("G9").Value = 13.75374933333333
$ws.Range("H9").Value = 41.261248
$ws.Range("I9").Value = 0.2533018145204882
$ws.Range("J9").Value = 0.2533018145204882
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 3527.070693515463
$ws.Range("R9").Value = 31743.63624163917
$ws.Range("S9").Value = 0.2113743497340395
$ws.Range("T9").Value = 0.2113743497340395
$ws.Range("G10").Value = 13.75374933333333
$ws.Range("H10").Value = 41.261248
$ws.Range("I10").Value = 0.2533018145204882
$ws.Range("J10").Value = 0.2533018145204882
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 671.7367713527397
$ws.Range("R10").Value = 6045.630942174657
$ws.Range("S10").Value = 0.04025661393693471
$ws.Range("T10").Value = 0.04025661393693469
